$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.00100255012512207
$ws.Range("C2").Value = 0.00199723243713378
$ws.Range("D2").Value = 0

$ws.Range("B3").Value = 0.00600075721740722
$ws.Range("C3").Value = 0.133010864257812
$ws.Range("D3").Value = 0.00100088119506835

$ws.Range("B4").Value = 0.0700054168701171
$ws.Range("C4").Value = 0.135011911392211
$ws.Range("D4").Value = 0.00900030136108398

$ws.Range("B5").Value = 0.902074337005615
$ws.Range("C5").Value = 0.202017545700073
$ws.Range("D5").Value = 0.0680086612701416

$ws.Range("J14").Select() | Out-Null
